# Fix the "2050" column label (previously holding a stray numeric value)
# on every table, and drop the "Total" row from each table.

$wb = $excel.ActiveWorkbook

# Helper: write a text label into a header cell without Excel
# auto-converting a numeric-looking string ("2050") into a number,
# and without changing the cell's style index (keep it identical to
# the bold/bordered header style already used by the row).
function Set-HeaderLabel {
    param($ws, [string]$cellAddr, [string]$formatSourceAddr, [string]$text)

    $ws.Range($cellAddr).NumberFormat = "@"
    $ws.Range($cellAddr).Value = $text
    $ws.Range($formatSourceAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)   # xlPasteFormats
}

# --- Sheets whose last header column (E1) changes "2040" -> "2050" ---
$simpleSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)
foreach ($name in $simpleSheets) {
    $ws = $wb.Worksheets.Item($name)
    Set-HeaderLabel $ws "E1" "D1" "2050"
}

# --- "Potencia Incremental" sheet uses period ranges: "2031-2040" -> "2041-2050" ---
$wsIncr = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-HeaderLabel $wsIncr "E1" "D1" "2041-2050"

# --- Remove the "Total" row (row 13) from the four tables that have one ---
$totalRow13Sheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)
foreach ($name in $totalRow13Sheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

# --- Remove the "Total" row (row 4) from the "Custo Total" table ---
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
